# Fix deleting active tahun-pelajaran issue:
# system should not delete an active tahun-pelajaran data.
#
# The previously-active class sheet ("12 MIA 2") is preserved as a brand
# new trailing sheet instead of being overwritten, while the other three
# class sheets are renamed/refreshed with the new tahun-ajaran roster.

$YELLOW = 65535      # Belum Lunas / Rp. 420.000 fill (FFFF00)
$BLUE   = 16764057   # Lunas fill (99CCFF)

function Set-StatusCell($cell, $text, $isLunas) {
    $cell.Value = $text
    if ($isLunas) {
        $cell.Interior.Color = $BLUE
    } else {
        $cell.Interior.Color = $YELLOW
    }
}

# $months is an array of 7 booleans (C..J = DaftarUlang, then Jan..Jun) saying
# whether that column should read "Lunas" (true) vs the default (false).
function Set-DataRow($ws, $rowNum, $no, $name, $months) {
    $ws.Cells.Item($rowNum, 1).Value = $no
    $ws.Cells.Item($rowNum, 2).Value = $name

    $cCell = $ws.Cells.Item($rowNum, 3)
    Set-StatusCell $cCell "Rp. 420.000" $months[0]
    if ($months[0]) { $cCell.Value = "Lunas" }

    $labels = @("JAN", "FEB", "MAR", "APR", "MEI", "JUN")
    for ($i = 1; $i -le 6; $i++) {
        $col = 3 + $i
        $cell = $ws.Cells.Item($rowNum, $col)
        $isLunas = $months[$i]
        if ($isLunas) {
            Set-StatusCell $cell "Lunas" $true
        } else {
            Set-StatusCell $cell "Belum Lunas" $false
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "27 jmy 6" -> "24 wao 6"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "24 wao 6"
$ws1.Range("A2").Value = "TAHUN AJARAN 2022/2023"
$ws1.Range("A5").Value = "KELAS 24 wao 6"

Set-DataRow $ws1 8  1 "Farah Prastuti"         @($false, $false, $false, $false, $false, $false, $false)
Set-DataRow $ws1 9  2 "Almira Padmasari S.E.I" @($false, $false, $false, $false, $false, $false, $false)
Set-DataRow $ws1 10 3 "Wahyu Dongoran"         @($false, $false, $false, $false, $false, $false, $false)
Set-DataRow $ws1 11 4 "Citra Yolanda"          @($false, $false, $false, $false, $false, $false, $false)

$ws1.Activate()
$ws1.Range("I11").Select()

# ---------------------------------------------------------------------
# Sheet 2: "12 MIA 2" -> "45 iie 1"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "45 iie 1"
$ws2.Range("A2").Value = "TAHUN AJARAN 2022/2023"
$ws2.Range("A5").Value = "KELAS 45 iie 1"

Set-DataRow $ws2 8  1 "Diana Padmasari"           @($true,  $false, $false, $false, $false, $false, $false)
Set-DataRow $ws2 9  2 "Saiful Sihotang"            @($false, $false, $false, $false, $false, $false, $false)
Set-DataRow $ws2 10 3 "Cakrabirawa Narpati S.Ked"  @($false, $false, $false, $false, $false, $false, $false)
Set-DataRow $ws2 11 4 "Kamila Aryani"              @($false, $false, $false, $false, $false, $false, $false)
Set-DataRow $ws2 12 5 "Zulaikha Zalindra Yuniar"   @($false, $false, $false, $false, $false, $false, $false)

$ws2.Activate()
$ws2.Range("I12").Select()

# ---------------------------------------------------------------------
# Sheet 3: "04 gbk 2" -> "49 sbj 7"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "49 sbj 7"
$ws3.Range("A2").Value = "TAHUN AJARAN 2022/2023"
$ws3.Range("A5").Value = "KELAS 49 sbj 7"

Set-DataRow $ws3 8  1 "Koko Wahyudin"                        @($false, $false, $false, $false, $false, $false, $false)
Set-DataRow $ws3 9  2 "Tasdik Lazuardi"                      @($false, $false, $false, $false, $false, $false, $false)
Set-DataRow $ws3 10 3 "Kenari Darsirah Situmorang M.Kom."    @($true,  $false, $false, $false, $false, $false, $false)
Set-DataRow $ws3 11 4 "Danang Wahyudin"                      @($false, $false, $false, $false, $false, $false, $false)
Set-DataRow $ws3 12 5 "Ibun Dongoran"                        @($false, $false, $false, $false, $false, $false, $false)
Set-DataRow $ws3 13 6 "Sarah Agustina"                       @($false, $false, $false, $false, $false, $false, $false)

$ws3.Activate()
$ws3.Range("I13").Select()

# ---------------------------------------------------------------------
# Sheet 4 (NEW): re-create "12 MIA 2" as its own sheet so the active
# tahun-pelajaran roster isn't lost when sheet 2 became "45 iie 1".
# ---------------------------------------------------------------------
$ws3.Copy($null, $ws3)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "12 MIA 2"
$ws4.Range("A2").Value = "TAHUN AJARAN 2022/2023"
$ws4.Range("A5").Value = "KELAS 12 MIA 2"

Set-DataRow $ws4 8  1 "Himawan Maheswara"    @($false, $false, $false, $false, $false, $false, $false)
Set-DataRow $ws4 9  2 "Dagel Tasnim Wasita"  @($false, $false, $false, $false, $false, $false, $false)
Set-DataRow $ws4 10 3 "Harjasa Mangunsong"   @($true,  $true,  $false, $false, $false, $false, $false)
Set-DataRow $ws4 11 4 "Nalar Lazuardi"       @($false, $false, $false, $false, $false, $false, $false)
Set-DataRow $ws4 12 5 "Maida Unjani Sudiati" @($false, $false, $false, $false, $false, $false, $false)
Set-DataRow $ws4 13 6 "Wd"                   @($false, $false, $false, $false, $false, $false, $false)

$ws4.Activate()
$ws4.Range("I13").Select()
